$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "44.028.64"
Set-TextCell $ws.Range("E2") "  -0.42%  "
Set-TextCell $ws.Range("D3") "2.238.27"
Set-TextCell $ws.Range("E3") "  -1.13%  "
Set-TextCell $ws.Range("E4") "  +0.08%  "
Set-TextCell $ws.Range("D5") "305.55"
Set-TextCell $ws.Range("E5") "  -4.69%  "
Set-TextCell $ws.Range("D6") "95.58"
Set-TextCell $ws.Range("E6") "  -6.94%  "
Set-TextCell $ws.Range("E7") "  -1.80%  "
Set-TextCell $ws.Range("E8") "  +0.21%  "
Set-TextCell $ws.Range("D9") "0.523"
Set-TextCell $ws.Range("E9") "  -5.65%  "
Set-TextCell $ws.Range("D10") "34.94"
Set-TextCell $ws.Range("E10") "  -6.69%  "
Set-TextCell $ws.Range("E11") "  -3.68%  "
Set-TextCell $ws.Range("D12") "7.21"
Set-TextCell $ws.Range("E12") "  -5.64%  "
Set-TextCell $ws.Range("E13") "  -2.93%  "
Set-TextCell $ws.Range("D14") "2.580.43"
Set-TextCell $ws.Range("D15") "2.241.09"
Set-TextCell $ws.Range("E15") "  -0.85%  "
Set-TextCell $ws.Range("E16") "  -4.61%  "
Set-TextCell $ws.Range("D17") "13.55"
Set-TextCell $ws.Range("E17") "  -6.53%  "
Set-TextCell $ws.Range("D18") "43.736.06"
Set-TextCell $ws.Range("E18") "  -0.76%  "
Set-TextCell $ws.Range("E19") "  -3.20%  "
Set-TextCell $ws.Range("D20") "12.27"
Set-TextCell $ws.Range("E20") "  -8.68%  "
Set-TextCell $ws.Range("E21") "  -5.07%  "
Set-TextCell $ws.Range("D22") "64.79"
Set-TextCell $ws.Range("E22") "  -1.65%  "
Set-TextCell $ws.Range("D23") "236.35"
Set-TextCell $ws.Range("E23") "  +0.15%  "
Set-TextCell $ws.Range("E24") "  -7.62%  "
Set-TextCell $ws.Range("E25") "  -7.63%  "
Set-TextCell $ws.Range("E26") "  +0.55%  "
Set-TextCell $ws.Range("D27") "9.95"
Set-TextCell $ws.Range("E27") "  -3.52%  "
Set-TextCell $ws.Range("D28") "38.11"
Set-TextCell $ws.Range("E28") "  -3.38%  "
Set-TextCell $ws.Range("D29") "2.14"
Set-TextCell $ws.Range("E29") "  -1.93%  "
Set-TextCell $ws.Range("D30") "5.95"
Set-TextCell $ws.Range("E30") "  -5.07%  "
Set-TextCell $ws.Range("D31") "20.01"
Set-TextCell $ws.Range("E31") "  -1.39%  "
Set-TextCell $ws.Range("D32") "155.08"
Set-TextCell $ws.Range("E32") "  -4.74%  "
Set-TextCell $ws.Range("E33") "  -5.49%  "
Set-TextCell $ws.Range("E34") "  +8.03%  "
Set-TextCell $ws.Range("D35") "2.62"
Set-TextCell $ws.Range("E35") "  -2.34%  "
Set-TextCell $ws.Range("D36") "0.108"
Set-TextCell $ws.Range("E36") "  -6.15%  "
Set-TextCell $ws.Range("E37") "  -1.18%  "
Set-TextCell $ws.Range("E38") "  -10.59%  "
Set-TextCell $ws.Range("D39") "15.25"
Set-TextCell $ws.Range("E39") "  -9.40%  "
Set-TextCell $ws.Range("D40") "3.36"
Set-TextCell $ws.Range("E40") "  -9.54%  "
Set-TextCell $ws.Range("D41") "3.82"
Set-TextCell $ws.Range("E41") "  -9.89%  "
Set-TextCell $ws.Range("E42") "  -5.18%  "
Set-TextCell $ws.Range("E43") "  +0.29%  "
Set-TextCell $ws.Range("D44") "1.736.45"
Set-TextCell $ws.Range("E44") "  -2.89%  "
Set-TextCell $ws.Range("D45") "85.32"
Set-TextCell $ws.Range("E45") "  +3.12%  "
Set-TextCell $ws.Range("E46") "  -6.16%  "
Set-TextCell $ws.Range("D47") "99.95"
Set-TextCell $ws.Range("E47") "  -4.98%  "
Set-TextCell $ws.Range("D48") "4.92"
Set-TextCell $ws.Range("E48") "  -5.92%  "
Set-TextCell $ws.Range("B49") "ordi"
Set-TextCell $ws.Range("C49") "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-TextCell $ws.Range("D49") "68.99"
Set-TextCell $ws.Range("E49") "  -8.11%  "
Set-TextCell $ws.Range("B50") "FraxShare"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws.Range("D50") "8.07"
Set-TextCell $ws.Range("E50") "  -4.13%  "
Set-TextCell $ws.Range("B51") "MultiversX"
Set-TextCell $ws.Range("C51") "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell $ws.Range("D51") "54.24"
Set-TextCell $ws.Range("E51") "  -7.53%  "

Write-Host "Applied" 87 "cell updates"
